# Daily attendance processing - 2026-01-10 14:03:05
# Re-order the "Recorded By" (column G) comma-separated list of authors
# so that the literal token "System" (exact case) is moved to the front
# of the list, leaving the relative order of the remaining tokens intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -notmatch ",") { continue }

    $parts = $val -split ",\s*"
    $idx = [Array]::IndexOf($parts, "System")

    if ($idx -ge 0) {
        $rest = @()
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($i -ne $idx) {
                $rest += $parts[$i]
            }
        }
        $newParts = @("System") + $rest
        $cell.Value2 = ($newParts -join ", ")
    }
}
